$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.617.02"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "2.397.25"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.42"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.98"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").Value = "2.402.04"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.00"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "2.830.09"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "60.310.84"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "2.411.99"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.06"
$ws.Range("E19").Value = "  +6.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.60"
$ws.Range("E20").Value = "  -0.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.62"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.75"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "558.94"
$ws.Range("E27").Value = "  -5.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.97"
$ws.Range("E28").Value = "  -5.94%  "

$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("E32").Value = "  -2.56%  "

$ws.Range("E33").Value = "  -2.80%  "

$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  +1.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.30"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.369"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.09"
$ws.Range("E41").Value = "  -1.43%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.73"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("E45").Value = "  +4.75%  "

$ws.Range("D46").Value = "0.0₆0273"
$ws.Range("E46").Value = "  -4.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.31"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.586"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.18"
$ws.Range("E51").Value = "  -2.55%  "
